$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New client record: Chris Lucas
$ws.Range("B3").Value = 1592999959
$ws.Range("C3").Value = "Chris"
$ws.Range("D3").Value = "Lucas"
$ws.Range("E3").Value = "62 West Wallaby Street"
$ws.Range("F3").Value = "Wigan"
$ws.Range("G3").Value = "W1 GAN"
$ws.Range("H3").Value = "United Kingdom"
$ws.Range("I3").Value = "+44 07987 654321"
$ws.Range("J3").Value = "chris.lucas@qa.com "
$ws.Range("K3").Value = "Tech"
$zwsp = [char]0x200B
$ws.Range("L3").Value = "DAILY" + $zwsp

# Column B best-fit width (matches Excel's computed best-fit of 11 characters)
$ws.Columns.Item(2).ColumnWidth = 10.166666666666666

# Update selection
$ws.Range("N8").Select()
